# ExcelStyle dataFormat 옵션 추가
# Rebuild the sheet grid: insert a new "값" column (B) holding the former
# A-column free text, push the old short labels into column A, and add a
# brand-new "cost" column (C) with a numeric money value + its own
# hyperlink row.
#
# Style plan:
#   - original A2:A4 carry the italic "value" style -> stays on B/C
#   - original B2:B4 carry the bold   "label" style -> moves to column A
# Donor styles are staged on scratch cells first so the real paste targets
# (which overlap the donors) don't clobber a style before it's been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Stage the two donor styles on scratch cells (untouched so far) ---
$ws.Range("A2").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Apply styles -------------------------------------------------------
# Italic "value" style -> columns B & C, rows 2-5
$ws.Range("E1").Copy()
$ws.Range("B2:C5").PasteSpecial(-4122)

# Bold "label" style -> column A, rows 2-5
$ws.Range("E2").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
$ws.Range("E1:E2").Clear()

# Header row style (s=1) already covers A1:B1; extend it to the new C1.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Values -------------------------------------------------------------
$ws.Range("A1").Value = "a푸"
$ws.Range("B1").Value = "b파"
$ws.Range("C1").Value = "cost"

$ws.Range("A2").Value = "가"
$ws.Range("B2").Value = "블라블라?블~라~블~라~"

$ws.Range("A3").Value = "나"
$ws.Range("B3").Value = ""

$ws.Range("A4").Value = "다"
$ws.Range("B4").Value = "https://www.google.com"

$ws.Range("A5").Value = "라"
$ws.Range("B5").Value = "https://www.naver.com"
$ws.Range("C5").Value = 100000

# ---- Hyperlinks ---------------------------------------------------------
# Old hyperlink lived on A4; it now belongs on B4 (same target), plus a
# brand-new one on B5 pointing at Naver.
$ws.Range("A4").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.google.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.naver.com")

# Hyperlinks.Add re-styles the target cell with the builtin "Hyperlink"
# style; restore the intended italic "value" style afterwards.
$ws.Range("C2").Copy()
$ws.Range("B4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
